# Insert a new bullet item right after the "Also manages state animations
# for tiles" list item, under the "Separate class for game table
# management:" heading, matching the ListParagraph/numPr bullet style of
# its neighbours.

$d = $word.ActiveDocument

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Also manages state animations for tiles") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph 'Also manages state animations for tiles'"
}

# Inserting a new paragraph right after the anchor's range copies the
# anchor paragraph's formatting (style + numbering), same as Word does
# when you press Enter at the end of a list item.
$anchor.Range.InsertParagraphAfter()

$newPara = $anchor.Next()
$newPara.Range.Text = "Management of worm ($([char]0x201C)virus$([char]0x201D)), power-ups, bombs"
